$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ParticipantsTab" Cypher query in B2 ---
# NOTE: single-quoted here-string (@' ... '@) so PowerShell does not treat
# the backticks around the RETURN column aliases as escape characters.
$query = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['BW']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@
$ws.Range("B2").Value = $query

# --- Row 2 grew taller to fit the longer (18-line) query text ---
$ws.Rows.Item(2).RowHeight = 279

# --- Selection / scroll position moved from B3 to B4 ---
$ws.Range("B4").Select() | Out-Null
